$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.734397888183594
$ws.Range("B1").Value = 1.886769533157349
$ws.Range("C1").Value = 5.11505126953125
$ws.Range("D1").Value = 1.240546584129333
$ws.Range("E1").Value = 0.6432112455368042
